# Apply scheduled market-data refresh to the per-sheet Leve profit tables.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for
# the rows whose upstream prices changed, sheet by sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 1559
$ws.Range("I125").Value = 754.875
$ws.Range("K125").Value = 6793.875
$ws.Range("M125").Value = -4333.875

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 584.4
$ws.Range("I2").Value = 530.5
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 530.5
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -417.5
$ws.Range("N2").Value = -1026
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 2445.3794
$ws.Range("I61").Value = 2194.25
$ws.Range("K61").Value = 2194.25
$ws.Range("M61").Value = -1982.25
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 3011.4167
$ws.Range("I74").Value = 2047.8334
$ws.Range("J74").Value = 3975
$ws.Range("K74").Value = 2047.8334
$ws.Range("L74").Value = 3975
$ws.Range("M74").Value = -1173.8334
$ws.Range("N74").Value = -5723
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 3011.4167
$ws.Range("I77").Value = 2047.8334
$ws.Range("J77").Value = 3975
$ws.Range("K77").Value = 10239.167
$ws.Range("L77").Value = 19875
$ws.Range("M77").Value = -5871.166999999999
$ws.Range("N77").Value = -28611
# Row 98 (Leve Item ID 18371)
$ws.Range("H98").Value = 63666.332
$ws.Range("J98").Value = 63666.332
$ws.Range("L98").Value = 63666.332
$ws.Range("N98").Value = -69656.33199999999
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 584.4
$ws.Range("I116").Value = 530.5
$ws.Range("J116").Value = 800
$ws.Range("K116").Value = 530.5
$ws.Range("L116").Value = 800
$ws.Range("M116").Value = 1763.5
$ws.Range("N116").Value = -5388
# Row 118 (Leve Item ID 26150)
$ws.Range("H118").Value = 45747.5
$ws.Range("J118").Value = 45747.5
$ws.Range("L118").Value = 45747.5
$ws.Range("N118").Value = -49061.5
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 2445.3794
$ws.Range("I136").Value = 2194.25
$ws.Range("K136").Value = 6582.75
$ws.Range("M136").Value = -4032.75

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 584.4
$ws.Range("I3").Value = 530.5
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 530.5
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = -416.5
$ws.Range("N3").Value = -1028
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1175386.2
$ws.Range("I99").Value = 39998.69
$ws.Range("J99").Value = 4127394
$ws.Range("K99").Value = 39998.69
$ws.Range("L99").Value = 4127394
$ws.Range("M99").Value = -38500.69
$ws.Range("N99").Value = -4130390
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 29791.916
$ws.Range("I105").Value = 51513.1
$ws.Range("K105").Value = 51513.1
$ws.Range("M105").Value = -49766.1
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 1509.8
$ws.Range("I107").Value = 1032.6666
$ws.Range("K107").Value = 1032.6666
$ws.Range("M107").Value = 887.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 28 (Leve Item ID 18348)
$ws.Range("H28").Value = 18969.8
$ws.Range("J28").Value = 18969.8
$ws.Range("L28").Value = 18969.8
$ws.Range("N28").Value = -19459.8
# Row 74 (Leve Item ID 10636)
$ws.Range("H74").Value = 61597
$ws.Range("J74").Value = 65916.39999999999
$ws.Range("L74").Value = 65916.39999999999
$ws.Range("N74").Value = -67664.39999999999
# Row 77 (Leve Item ID 10636)
$ws.Range("H77").Value = 61597
$ws.Range("J77").Value = 65916.39999999999
$ws.Range("L77").Value = 197749.2
$ws.Range("N77").Value = -206485.2
# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 5997.222
$ws.Range("I86").Value = 5233.222
$ws.Range("K86").Value = 5233.222
$ws.Range("M86").Value = -4110.222
# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 5997.222
$ws.Range("I89").Value = 5233.222
$ws.Range("K89").Value = 26166.11
$ws.Range("M89").Value = -20550.11
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2213.3125
$ws.Range("I132").Value = 2213.3125
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6639.9375
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4109.9375
$ws.Range("N132").Value = $null
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 2757.205
$ws.Range("I134").Value = 2757.394
$ws.Range("K134").Value = 8272.181999999999
$ws.Range("M134").Value = -5737.181999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 17 (Leve Item ID 4640)
$ws.Range("H17").Value = 1541
$ws.Range("I17").Value = 623
$ws.Range("K17").Value = 1869
$ws.Range("M17").Value = -1700
# Row 112 (Leve Item ID 27855)
$ws.Range("H112").Value = 4999.846
$ws.Range("I112").Value = 2999.3333
$ws.Range("K112").Value = 8997.999899999999
$ws.Range("M112").Value = -7889.999899999999
# Row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 9372
$ws.Range("J114").Value = 20717.75
$ws.Range("L114").Value = 62153.25
$ws.Range("N114").Value = -68661.25
# Row 117 (Leve Item ID 27870)
$ws.Range("H117").Value = 1244.8
$ws.Range("J117").Value = 1728.4
$ws.Range("L117").Value = 5185.200000000001
$ws.Range("N117").Value = -12069.2

$ws = $wb.Worksheets.Item("GSM")
# Row 10 (Leve Item ID 4306)
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 10000
$ws.Range("K10").Value = 10000
$ws.Range("M10").Value = -9831
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 2813.2666
$ws.Range("I102").Value = 2823.8462
$ws.Range("J102").Value = 2744.5
$ws.Range("K102").Value = 2823.8462
$ws.Range("L102").Value = 2744.5
$ws.Range("M102").Value = -1201.8462
$ws.Range("N102").Value = -5988.5
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 4398.364
$ws.Range("I122").Value = 3470.8333
$ws.Range("J122").Value = 8572.25
$ws.Range("K122").Value = 10412.4999
$ws.Range("L122").Value = 25716.75
$ws.Range("M122").Value = -7962.499899999999
$ws.Range("N122").Value = -30616.75
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 6372.4287
$ws.Range("I132").Value = 5966.6665
$ws.Range("J132").Value = 6676.75
$ws.Range("K132").Value = 17899.9995
$ws.Range("L132").Value = 20030.25
$ws.Range("M132").Value = -15369.9995
$ws.Range("N132").Value = -25090.25

$ws = $wb.Worksheets.Item("LTW")
# Row 74 (Leve Item ID 11990)
$ws.Range("H74").Value = 63333
$ws.Range("J74").Value = 65999.8
$ws.Range("L74").Value = 65999.8
$ws.Range("N74").Value = -67995.8
# Row 77 (Leve Item ID 11990)
$ws.Range("H77").Value = 63333
$ws.Range("J77").Value = 65999.8
$ws.Range("L77").Value = 197999.4
$ws.Range("N77").Value = -207983.4
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3058.6155
$ws.Range("I132").Value = 1782.875
$ws.Range("K132").Value = 5348.625
$ws.Range("M132").Value = -2818.625
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 2362.4348
$ws.Range("I136").Value = 2033.5385
$ws.Range("K136").Value = 6100.6155
$ws.Range("M136").Value = -3550.6155

$ws = $wb.Worksheets.Item("WVR")
# Row 93 (Leve Item ID 19613)
$ws.Range("H93").Value = 78784.664
$ws.Range("I93").Value = 73177
$ws.Range("J93").Value = 90000
$ws.Range("K93").Value = 73177
$ws.Range("L93").Value = 90000
$ws.Range("M93").Value = -70681
$ws.Range("N93").Value = -94992
# Row 95 (Leve Item ID 18243)
$ws.Range("H95").Value = 49985.715
$ws.Range("J95").Value = 49985.715
$ws.Range("L95").Value = 49985.715
$ws.Range("N95").Value = -55477.715
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 2613.05
$ws.Range("I122").Value = 1911.6428
$ws.Range("K122").Value = 5734.928400000001
$ws.Range("M122").Value = -3284.928400000001
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1146039.5
$ws.Range("I132").Value = 1527.08
$ws.Range("K132").Value = 4581.24
$ws.Range("M132").Value = -2051.24
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1353.5
$ws.Range("I136").Value = 791.2632
$ws.Range("K136").Value = 2373.7896
$ws.Range("M136").Value = 176.2103999999999
